$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 298010.25
$ws.Range("J17").Value = 298010.25
$ws.Range("L17").Value = 894030.75
$ws.Range("N17").Value = -894366.75
$ws.Range("H62").Value = 9344.725
$ws.Range("I62").Value = 10175.971
$ws.Range("J62").Value = 4634.3335
$ws.Range("K62").Value = 10175.971
$ws.Range("L62").Value = 4634.3335
$ws.Range("M62").Value = -9551.971
$ws.Range("N62").Value = -5882.3335
$ws.Range("H65").Value = 9344.725
$ws.Range("I65").Value = 10175.971
$ws.Range("J65").Value = 4634.3335
$ws.Range("K65").Value = 50879.855
$ws.Range("L65").Value = 23171.6675
$ws.Range("M65").Value = -47759.855
$ws.Range("N65").Value = -29411.6675
$ws.Range("H137").Value = 29948.947
$ws.Range("I137").Value = 50945.2
$ws.Range("J137").Value = 6619.778
$ws.Range("K137").Value = 152835.6
$ws.Range("L137").Value = 19859.334
$ws.Range("M137").Value = -150285.6
$ws.Range("N137").Value = -24959.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1074.2
$ws.Range("I32").Value = 733.38635
$ws.Range("J32").Value = 3573.5
$ws.Range("K32").Value = 733.38635
$ws.Range("L32").Value = 3573.5
$ws.Range("M32").Value = -446.38635
$ws.Range("N32").Value = -4147.5
$ws.Range("H76").Value = 39912.668
$ws.Range("J76").Value = 39912.668
$ws.Range("L76").Value = 39912.668
$ws.Range("N76").Value = -40588.668
$ws.Range("H79").Value = 39912.668
$ws.Range("J79").Value = 39912.668
$ws.Range("L79").Value = 39912.668
$ws.Range("N79").Value = -42252.668
$ws.Range("H132").Value = 1290494.9
$ws.Range("I132").Value = 1596196
$ws.Range("J132").Value = 439848.34
$ws.Range("K132").Value = 4788588
$ws.Range("L132").Value = 1319545.02
$ws.Range("M132").Value = -4786058
$ws.Range("N132").Value = -1324605.02

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1219.76
$ws.Range("I20").Value = 1302.0588
$ws.Range("J20").Value = 1044.875
$ws.Range("K20").Value = 1302.0588
$ws.Range("L20").Value = 1044.875
$ws.Range("M20").Value = -1055.0588
$ws.Range("N20").Value = -1538.875
$ws.Range("H94").Value = 1564
$ws.Range("I94").Value = 555.5833
$ws.Range("J94").Value = 3984.2
$ws.Range("K94").Value = 555.5833
$ws.Range("L94").Value = 3984.2
$ws.Range("M94").Value = -104.5833
$ws.Range("N94").Value = -4886.2
$ws.Range("H134").Value = 22036.133
$ws.Range("I134").Value = 1013.175
$ws.Range("J134").Value = 86722.16
$ws.Range("K134").Value = 3039.525
$ws.Range("L134").Value = 260166.48
$ws.Range("M134").Value = -504.5249999999996
$ws.Range("N134").Value = -265236.48

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8973.547
$ws.Range("I31").Value = 7602.2705
$ws.Range("J31").Value = 12144.625
$ws.Range("K31").Value = 7602.2705
$ws.Range("L31").Value = 12144.625
$ws.Range("M31").Value = -7307.2705
$ws.Range("N31").Value = -12734.625
$ws.Range("H34").Value = 8973.547
$ws.Range("I34").Value = 7602.2705
$ws.Range("J34").Value = 12144.625
$ws.Range("K34").Value = 7602.2705
$ws.Range("L34").Value = 12144.625
$ws.Range("M34").Value = -7400.2705
$ws.Range("N34").Value = -12548.625
$ws.Range("H132").Value = 1004.8542
$ws.Range("I132").Value = 1005.4722
$ws.Range("K132").Value = 3016.4166
$ws.Range("M132").Value = -486.4166

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7086.826
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 10230.538
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 10230.538
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -12226.538
$ws.Range("H83").Value = 7086.826
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 10230.538
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 51152.69
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -61136.69
$ws.Range("H102").Value = 27082.354
$ws.Range("I102").Value = 11952.77
$ws.Range("J102").Value = 76253.5
$ws.Range("K102").Value = 11952.77
$ws.Range("L102").Value = 76253.5
$ws.Range("M102").Value = -10330.77
$ws.Range("N102").Value = -79497.5
$ws.Range("H122").Value = 1194.0769
$ws.Range("I122").Value = 999.8570999999999
$ws.Range("J122").Value = 1420.6666
$ws.Range("K122").Value = 2999.5713
$ws.Range("L122").Value = 4261.9998
$ws.Range("M122").Value = -549.5712999999996
$ws.Range("N122").Value = -9161.9998
$ws.Range("H132").Value = 23746.29
$ws.Range("I132").Value = 1470.3077
$ws.Range("J132").Value = 54229.21
$ws.Range("K132").Value = 4410.9231
$ws.Range("L132").Value = 162687.63
$ws.Range("M132").Value = -1880.9231
$ws.Range("N132").Value = -167747.63

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1961.6538
$ws.Range("I68").Value = 1753.3334
$ws.Range("J68").Value = 2245.7273
$ws.Range("K68").Value = 1753.3334
$ws.Range("L68").Value = 2245.7273
$ws.Range("M68").Value = -1004.3334
$ws.Range("N68").Value = -3743.7273
$ws.Range("H71").Value = 1961.6538
$ws.Range("I71").Value = 1753.3334
$ws.Range("J71").Value = 2245.7273
$ws.Range("K71").Value = 8766.666999999999
$ws.Range("L71").Value = 11228.6365
$ws.Range("M71").Value = -5022.666999999999
$ws.Range("N71").Value = -18716.6365
$ws.Range("H122").Value = 2639.2195
$ws.Range("I122").Value = 2666.7646
$ws.Range("J122").Value = 2505.4285
$ws.Range("K122").Value = 8000.293799999999
$ws.Range("L122").Value = 7516.2855
$ws.Range("M122").Value = -5550.293799999999
$ws.Range("N122").Value = -12416.2855
$ws.Range("H136").Value = 141762.92
$ws.Range("I136").Value = 189286.1
$ws.Range("J136").Value = 1833.5555
$ws.Range("K136").Value = 567858.3
$ws.Range("L136").Value = 5500.666499999999
$ws.Range("M136").Value = -565308.3
$ws.Range("N136").Value = -10600.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2602.22
$ws.Range("I132").Value = 407.38235
$ws.Range("J132").Value = 7266.25
$ws.Range("K132").Value = 1222.14705
$ws.Range("L132").Value = 21798.75
$ws.Range("M132").Value = 1307.85295
$ws.Range("N132").Value = -26858.75
$ws.Range("H136").Value = 857854.5600000001
$ws.Range("I136").Value = 1066780.6
$ws.Range("J136").Value = 357924.44
$ws.Range("K136").Value = 3200341.8
$ws.Range("L136").Value = 1073773.32
$ws.Range("M136").Value = -3197791.8
$ws.Range("N136").Value = -1078873.32
